$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused trailing columns (L:Q) so the sheet shrinks from A1:Q6 to A1:K6
$ws.Range("L1:Q6").EntireColumn.Delete()

# Row 1
$ws.Range("A1").Value = "climate_zone"
$ws.Range("B1").Value = "iso"
$ws.Range("C1").Value = "CSCC"
$ws.Range("D1").Value = "SLAND"
$ws.Range("E1").Value = "ELUC"
$ws.Range("F1").Value = "EFOS"
$ws.Range("G1").Value = "F_ab"
$ws.Range("H1").Value = "F_ac"
$ws.Range("I1").Value = "F_abc"
$ws.Range("J1").Value = "continent"
$ws.Range("K1").Value = "economic_group"

# Row 2
$ws.Range("A2").Value = "'1"
$ws.Range("B2").Value = "AGOARGBDIBENBFABGDBHSBLZBOLBRABRNCAFCIVCMRCODCOGCOKCOLCOMCRICUBDOMECUETHFJIGABGHAGINGMBGNBGNQGTMGUYHNDHTIIDNINDJAMKENKHMKIRLAOLBRLKAMDGMEXMMRMOZMUSMWIMYSNGANICPANPERPHLPNGPRYRWASENSLBSLESLVSTPSURSWZTGOTHATLSTTOTZAUGAURYVCTVENVNMVUTWSMZMB"
$ws.Range("C2").Value = 479.5156661484671
$ws.Range("D2").Value = 1.831527709960938
$ws.Range("E2").Value = -1.11149014075
$ws.Range("F2").Value = -1.641572104287398
$ws.Range("G2").Value = 0.7200375850197459
$ws.Range("H2").Value = -2.753062245037398
$ws.Range("I2").Value = -0.9215345192676525
$ws.Range("J2").Value = "AfricaLatin America and the CaribbeanAfricaAfricaAfricaAsiaLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanAsiaAfricaAfricaAfricaAfricaAfricaOceaniaLatin America and the CaribbeanAfricaLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanAfricaOceaniaAfricaAfricaAfricaAfricaAfricaAfricaLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanAsiaAsiaLatin America and the CaribbeanAfricaAsiaOceaniaAsiaAfricaAsiaAfricaLatin America and the CaribbeanAsiaAfricaAfricaAfricaAsiaAfricaLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanAsiaOceaniaLatin America and the CaribbeanAfricaAfricaOceaniaAfricaLatin America and the CaribbeanAfricaLatin America and the CaribbeanAfricaAfricaAsiaAsiaLatin America and the CaribbeanAfricaAfricaLatin America and the CaribbeanLatin America and the CaribbeanLatin America and the CaribbeanAsiaOceaniaOceaniaAfrica"
$ws.Range("K2").Value = "LDCOtherLDCLDCLDCOtherOtherOtherOtherBRICSOtherLDCOtherOtherLDCOtherOtherOECDLDCOECDOtherOtherOtherBRICSOtherOtherOtherLDCLDCLDCOtherOtherOtherOtherLDCOtherBRICSOtherOtherLDCOtherLDCLDCOtherLDCOECDLDCLDCOtherLDCOtherOtherOtherOtherOtherOtherOtherOtherLDCLDCOtherLDCOtherLDCOtherLDCLDCOtherLDCOtherLDCLDCOtherOtherOtherOtherOtherOtherLDC"

# Row 3
$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "AREAUSBWACHLCPVDJIDZAEGYERIIRNIRQISRJORKWTLBYLSOMARMLIMRTNAMNEROMNPAKPSEQATSAUSDNSOMSOMSSDSYRTCDTKMTUNYEMZAFZWE"
$ws.Range("C3").Value = 134.9298307888636
$ws.Range("D3").Value = 0.1445359289646149
$ws.Range("E3").Value = -0.05673168899999999
$ws.Range("F3").Value = -1.040377085387413
$ws.Range("G3").Value = 0.0878042482414714
$ws.Range("H3").Value = -1.097108774387412
$ws.Range("I3").Value = -0.9525728371459411
$ws.Range("J3").Value = "AsiaOceaniaAfricaLatin America and the CaribbeanAfricaAfricaAfricaAfricaAfricaAsiaAsiaAsiaAsiaAsiaAfricaAfricaAfricaAfricaAfricaAfricaAfricaAsiaAsiaAsiaAsiaAsiaAfricaAfricaAfricaAfricaAsiaAfricaAsiaAfricaAsiaAfricaAfrica"
$ws.Range("K3").Value = "BRICSOECDOtherOECDOtherLDCOtherBRICSLDCBRICSOtherOECDOtherOtherOtherLDCOtherLDCOtherOtherLDCOtherOtherOtherOtherOtherOtherOtherOtherLDCOtherLDCOtherOtherOtherBRICSLDC"

# Row 4
$ws.Range("A4").Value = "'3"
$ws.Range("B4").Value = "ALBAUTBELBGRBIHCHECYPCZEDEUDNKESPFRAGBRGRCHRVHUNIRLITAJPNKORLBNLUXMKDMNENLDNZLPRTROUSRBSVKSVNUSA"
$ws.Range("C4").Value = 25.15215896558901
$ws.Range("D4").Value = 0.3941102623939514
$ws.Range("E4").Value = 0.04610860625
$ws.Range("F4").Value = -2.757478497830251
$ws.Range("G4").Value = 0.4402188742900945
$ws.Range("H4").Value = -2.711369891580251
$ws.Range("I4").Value = -2.317259623540156
$ws.Range("J4").Value = "EuropeEuropeEuropeEuropeEuropeEuropeAsiaEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeEuropeAsiaAsiaAsiaEuropeEuropeEuropeEuropeOceaniaEuropeEuropeEuropeEuropeEuropeNorth America"
$ws.Range("K4").Value = "OtherOECDOECDOtherOtherOECDOtherOECDOECDOECDOECDOECDOECDOECDOtherOECDOECDOECDOECDOECDOtherOECDOtherOtherOECDOECDOECDOtherOtherOECDOECDOECD"

# Row 5
$ws.Range("A5").Value = "'4"
$ws.Range("B5").Value = "AFGARMAZEBLRBTNCHNESTFINGEOKAZKGZLTULVAMDAMNGNPLPOLPRKSWETJKTURUKRUZB"
$ws.Range("C5").Value = 54.78011043896012
$ws.Range("D5").Value = 0.3781403601169586
$ws.Range("E5").Value = 0.03537037825
$ws.Range("F5").Value = -3.310600968240123
$ws.Range("G5").Value = 0.4135107224908191
$ws.Range("H5").Value = -3.275230589990122
$ws.Range("I5").Value = -2.897090245749303
$ws.Range("J5").Value = "AsiaAsiaAsiaEuropeAsiaAsiaEuropeEuropeAsiaAsiaAsiaEuropeEuropeEuropeAsiaAsiaEuropeAsiaEuropeAsiaAsiaEuropeAsia"
$ws.Range("K5").Value = "LDCOtherOtherOtherOtherBRICSOECDOECDOtherOtherOtherOECDOECDOtherOtherOtherOECDOtherOECDOtherOECDOtherOther"

# Row 6
$ws.Range("A6").Value = "'5"
$ws.Range("B6").Value = "CANISLNORRUS"
$ws.Range("C6").Value = -31.93801035976101
$ws.Range("D6").Value = 0.6215693950653076
$ws.Range("E6").Value = -0.0877665435
$ws.Range("F6").Value = -0.6239657808630126
$ws.Range("G6").Value = 0.5338028868391501
$ws.Range("H6").Value = -0.7117323243630125
$ws.Range("I6").Value = -0.09016289402386246
$ws.Range("J6").Value = "North AmericaEuropeEuropeAsia"
$ws.Range("K6").Value = "OECDOECDOECDBRICS"
